# Weekly fruit/vegetable price update:
# Insert two new report rows (week's data) above the existing row 39,
# shifting the rest of the table down by two rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A39:A40").EntireRow.Insert()

# New row 39
$ws.Cells.Item(39, 1).Value = 5
$ws.Cells.Item(39, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(39, 3).Value = "Maule"
$ws.Cells.Item(39, 4).Value2 = 45027
$ws.Cells.Item(39, 5).Value = 7
$ws.Cells.Item(39, 6).Value = 100112043
$ws.Cells.Item(39, 7).Value = "Pepino dulce"
$ws.Cells.Item(39, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(39, 9).Value = "Primera"
$ws.Cells.Item(39, 10).Value = 400
$ws.Cells.Item(39, 11).Value = 13000
$ws.Cells.Item(39, 12).Value = 13000
$ws.Cells.Item(39, 13).Value = 13000
$ws.Cells.Item(39, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(39, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(39, 16).Value = 722
$ws.Cells.Item(39, 17).Value = 18
$ws.Cells.Item(39, 18).Value = "Hortaliza"

# New row 40
$ws.Cells.Item(40, 1).Value = 5
$ws.Cells.Item(40, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(40, 3).Value = "Maule"
$ws.Cells.Item(40, 4).Value2 = 45027
$ws.Cells.Item(40, 5).Value = 7
$ws.Cells.Item(40, 6).Value = 100112043
$ws.Cells.Item(40, 7).Value = "Pepino dulce"
$ws.Cells.Item(40, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(40, 9).Value = "Segunda"
$ws.Cells.Item(40, 10).Value = 200
$ws.Cells.Item(40, 11).Value = 11000
$ws.Cells.Item(40, 12).Value = 11000
$ws.Cells.Item(40, 13).Value = 11000
$ws.Cells.Item(40, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(40, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(40, 16).Value = 611
$ws.Cells.Item(40, 17).Value = 18
$ws.Cells.Item(40, 18).Value = "Hortaliza"

"done"
